# Uniformize valency-pattern labels (column L) and fill in the matching
# X/Y gloss columns (I/J) that were previously left blank, per the commit:
# "All structural tables have been made uniform in terms of how valency
# classes, X and Y columns etc. are organized."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old compact valency-pattern label -> new underscored label.
$labelMap = @{
    "POSSNOM"  = "POSS_NOM"
    "NOMINSTR" = "NOM_INS"
    "NOMDIR"   = "NOM_DIR"
    "NOMABL"   = "NOM_ABL"
    "NOMLOC"   = "NOM_LOC"
    "NOMORN"   = "NOM_ORN"
    "ACCNOM"   = "ACC_NOM"
    "DATNOM"   = "DAT_NOM"
    "NOMDAT"   = "NOM_DAT"
    "NOMNOM"   = "NOM_NOM"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellI = $ws.Cells.Item($r, 9)   # I = X
    $cellJ = $ws.Cells.Item($r, 10)  # J = Y
    $cellK = $ws.Cells.Item($r, 11)  # K = locus
    $cellL = $ws.Cells.Item($r, 12)  # L = valency_pattern

    $iVal = $cellI.Value2
    $jVal = $cellJ.Value2
    $lVal = $cellL.Value2

    if ($labelMap.ContainsKey($lVal)) {
        # Simple relabeling rows: I/J/K stay as-is, only L is renamed.
        $cellL.Value = $labelMap[$lVal]
    }
    elseif (($lVal -eq "TR") -and ($iVal -eq "TR") -and ($jVal -eq $null)) {
        # Transitive-verb rows: X/Y were left blank; fill in the
        # canonical NOM/ACC glosses (locus/valency_pattern stay "TR").
        $cellI.Value = "NOM"
        $cellJ.Value = "ACC"
    }
    elseif (($iVal -eq "*") -and ($jVal -eq $null) -and ($cellK.Value2 -eq "*")) {
        # "No valency info" rows: Y was left blank; mirror X/locus.
        $cellJ.Value = "*"
    }
}

# Reset the view: scroll back to A1 and clear any saved selection, since
# the sheet no longer starts scrolled to column I with N15 selected.
$ws.Activate()
$ws.Range("A1").Select()
